$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 44.09582266666666
$ws.Range("H2").Value = 132.287468
$ws.Range("I2").Value = 0.1927468402671175
$ws.Range("J2").Value = 0.1927468402671175
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.6847623333333334
$ws.Range("N2").Value = 2.054287
$ws.Range("O2").Value = 0.04097600788991114
$ws.Range("P2").Value = 0.04097600788991113
$ws.Range("Q2").Value = 30.19515841947955
$ws.Range("R2").Value = 271.756425775316
$ws.Range("S2").Value = 0.007897996047540851
$ws.Range("T2").Value = 0.007897996047540849

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 44.09582266666666
$ws.Range("H3").Value = 132.287468
$ws.Range("I3").Value = 0.1927468402671175
$ws.Range("J3").Value = 0.1927468402671175
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 14.24499333333333
$ws.Range("N3").Value = 42.73498
$ws.Range("O3").Value = 0.8524168617409322
$ws.Range("P3").Value = 0.8524168617409322
$ws.Range("Q3").Value = 628.1446999145155
$ws.Range("R3").Value = 5653.30229923064
$ws.Range("S3").Value = 0.1643006566909771
$ws.Range("T3").Value = 0.1643006566909771

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 44.09582266666666
$ws.Range("H4").Value = 132.287468
$ws.Range("I4").Value = 0.1927468402671175
$ws.Range("J4").Value = 0.1927468402671175
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.781543666666667
$ws.Range("N4").Value = 5.344631000000001
$ws.Range("O4").Value = 0.1066071303691566
$ws.Range("P4").Value = 0.1066071303691566
$ws.Range("Q4").Value = 78.55863359825645
$ws.Range("R4").Value = 707.027702384308
$ws.Range("S4").Value = 0.02054818752859961
$ws.Range("T4").Value = 0.0205481875285996

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 111.3149186666667
$ws.Range("H5").Value = 333.944756
$ws.Range("I5").Value = 0.4865676055026886
$ws.Range("J5").Value = 0.4865676055026886
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.6847623333333334
$ws.Range("N5").Value = 2.054287
$ws.Range("O5").Value = 0.04097600788991114
$ws.Range("P5").Value = 0.04097600788991113
$ws.Range("Q5").Value = 76.22426344099688
$ws.Range("R5").Value = 686.0183709689719
$ws.Range("S5").Value = 0.01993759804205334
$ws.Range("T5").Value = 0.01993759804205333

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 111.3149186666667
$ws.Range("H6").Value = 333.944756
$ws.Range("I6").Value = 0.4865676055026886
$ws.Range("J6").Value = 0.4865676055026886
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 14.24499333333333
$ws.Range("N6").Value = 42.73498
$ws.Range("O6").Value = 0.8524168617409322
$ws.Range("P6").Value = 0.8524168617409322
$ws.Range("Q6").Value = 1585.680274307209
$ws.Range("R6").Value = 14271.12246876488
$ws.Range("S6").Value = 0.4147584313074018
$ws.Range("T6").Value = 0.4147584313074018

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 111.3149186666667
$ws.Range("H7").Value = 333.944756
$ws.Range("I7").Value = 0.4865676055026886
$ws.Range("J7").Value = 0.4865676055026886
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.781543666666667
$ws.Range("N7").Value = 5.344631000000001
$ws.Range("O7").Value = 0.1066071303691566
$ws.Range("P7").Value = 0.1066071303691566
$ws.Range("Q7").Value = 198.3123883561151
$ws.Range("R7").Value = 1784.811495205036
$ws.Range("S7").Value = 0.05187157615323349
$ws.Range("T7").Value = 0.05187157615323348

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 40.34450033333334
$ws.Range("H8").Value = 121.033501
$ws.Range("I8").Value = 0.176349470111689
$ws.Range("J8").Value = 0.176349470111689
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.6847623333333334
$ws.Range("N8").Value = 2.054287
$ws.Range("O8").Value = 0.04097600788991114
$ws.Range("P8").Value = 0.04097600788991113
$ws.Range("Q8").Value = 27.62639418542078
$ws.Range("R8").Value = 248.637547668787
$ws.Range("S8").Value = 0.007226097278678216
$ws.Range("T8").Value = 0.007226097278678215

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 40.34450033333334
$ws.Range("H9").Value = 121.033501
$ws.Range("I9").Value = 0.176349470111689
$ws.Range("J9").Value = 0.176349470111689
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 14.24499333333333
$ws.Range("N9").Value = 42.73498
$ws.Range("O9").Value = 0.8524168617409322
$ws.Range("P9").Value = 0.8524168617409322
$ws.Range("Q9").Value = 574.7071382849978
$ws.Range("R9").Value = 5172.36424456498
$ws.Range("S9").Value = 0.1503232618822822
$ws.Range("T9").Value = 0.1503232618822822

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 40.34450033333334
$ws.Range("H10").Value = 121.033501
$ws.Range("I10").Value = 0.176349470111689
$ws.Range("J10").Value = 0.176349470111689
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.781543666666667
$ws.Range("N10").Value = 5.344631000000001
$ws.Range("O10").Value = 0.1066071303691566
$ws.Range("P10").Value = 0.1066071303691566
$ws.Range("Q10").Value = 71.87548905368124
$ws.Range("R10").Value = 646.8794014831311
$ws.Range("S10").Value = 0.01880011095072852
$ws.Range("T10").Value = 0.01880011095072852

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 33.02061066666666
$ws.Range("H11").Value = 99.061832
$ws.Range("I11").Value = 0.1443360841185049
$ws.Range("J11").Value = 0.144336084118505
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.6847623333333334
$ws.Range("N11").Value = 2.054287
$ws.Range("O11").Value = 0.04097600788991114
$ws.Range("P11").Value = 0.04097600788991113
$ws.Range("Q11").Value = 22.61127040819822
$ws.Range("R11").Value = 203.501433673784
$ws.Range("S11").Value = 0.005914316521638736
$ws.Range("T11").Value = 0.005914316521638736

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 33.02061066666666
$ws.Range("H12").Value = 99.061832
$ws.Range("I12").Value = 0.1443360841185049
$ws.Range("J12").Value = 0.144336084118505
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 14.24499333333333
$ws.Range("N12").Value = 42.73498
$ws.Range("O12").Value = 0.8524168617409322
$ws.Range("P12").Value = 0.8524168617409322
$ws.Range("Q12").Value = 470.3783788092622
$ws.Range("R12").Value = 4233.405409283359
$ws.Range("S12").Value = 0.1230345118602712
$ws.Range("T12").Value = 0.1230345118602712

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 33.02061066666666
$ws.Range("H13").Value = 99.061832
$ws.Range("I13").Value = 0.1443360841185049
$ws.Range("J13").Value = 0.144336084118505
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.781543666666667
$ws.Range("N13").Value = 5.344631000000001
$ws.Range("O13").Value = 0.1066071303691566
$ws.Range("P13").Value = 0.1066071303691566
$ws.Range("Q13").Value = 58.82765980266578
$ws.Range("R13").Value = 529.4489382239921
$ws.Range("S13").Value = 0.01538725573659501
$ws.Range("T13").Value = 0.01538725573659501

